# daily auto push: 2026-01-24 13:39 UTC
#
# A new reading was recorded for 2026/01/24 (Saturday) at 19:00, so a
# row is inserted just before the existing "2026/12/29" block (current
# row 689), shifting that block and everything after it down by one
# row. The new row holds: 2026/01/24, 土, 19, 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 689; rows 689:730 become 690:731
$ws.Rows.Item(689).Insert()

# Force column A to text so the date-like string "2026/01/24" is stored
# literally instead of being auto-converted to a date serial number.
$ws.Cells.Item(689, 1).NumberFormat = "@"
$ws.Cells.Item(689, 1).Value = "2026/01/24"
$ws.Cells.Item(689, 2).Value = "土"
$ws.Cells.Item(689, 3).Value = 19
$ws.Cells.Item(689, 4).Value = 17

# Drop back to the default "Normal" style so the inserted row doesn't
# retain a lingering text-number-format style the cell picked up above,
# matching the unstyled look of the other data rows.
$ws.Cells.Item(689, 1).Style = "Normal"
